$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheet 2: Ebay -> MortgageCalculator ---
$ws2.Name = "MortgageCalculator"

# --- Header row (row 1): plain text, default style ---
$ws2.Range("A1").Value = "purchase price"
$ws2.Range("B1").Value = "down payment%"
$ws2.Range("C1").Value = "Loan Type"
$ws2.Range("D1").Value = "Expected Mortgage Payment"

# --- Column A (purchase price), quote-prefixed numbers stored as text ---
$ws2.Range("A2").Value = "'300000"
$ws2.Range("A3").Value = "'300000"
$ws2.Range("A4").Value = "'400000"
$ws2.Range("A5").Value = "'400000"
$ws2.Range("A8").Value = "'325000"
$ws2.Range("A9").Value = "'325000"

# --- Column B (down payment %), quote-prefixed numbers stored as text ---
$ws2.Range("B2").Value = "'10"
$ws2.Range("B3").Value = "'20"
$ws2.Range("B4").Value = "'15"
$ws2.Range("B5").Value = "'20"
$ws2.Range("B6").Value = "'13"
$ws2.Range("B7").Value = "'10"
$ws2.Range("B8").Value = "'20"
$ws2.Range("B9").Value = "'20"

# --- Column A rows 6/7 (FHA): quote-prefixed number with explicit black Calibri font ---
$ws2.Range("A6").Value = "'250000"
$ws2.Range("A7").Value = "'250000"
$ws2.Range("A6:A7").Font.Color = 0

# --- Column D (expected payment), quote-prefixed text values with Text ("@") number format ---
$ws2.Range("D2").Value = "'`$1,490"
$ws2.Range("D3").Value = "'`$1,551"
$ws2.Range("D4").Value = "'`$2,742"
$ws2.Range("D5").Value = "'`$3,238"
$ws2.Range("D6").Value = "'`$915"
$ws2.Range("D7").Value = "'`$1,565"
$ws2.Range("D8").Value = "'`$1,266"
$ws2.Range("D9").Value = "'`$2,015"
$ws2.Range("D2:D9").NumberFormat = "@"

# --- Column C (loan type), plain text but with Menlo font matching DropDownList sheet ---
$ws2.Range("C2").Value = "30-Year Fixed"
$ws2.Range("C3").Value = "20-Year Fixed"
$ws2.Range("C4").Value = "15-Year Fixed"
$ws2.Range("C5").Value = "10-Year Fixed"
$ws2.Range("C6").Value = "FHA 30-Year Fixed"
$ws2.Range("C7").Value = "FHA 15-Year Fixed"
$ws2.Range("C8").Value = "VA 30-Year Fixed"
$ws2.Range("C9").Value = "VA 15-Year Fixed"
$ws1.Range("A2").Copy()
$ws2.Range("C2:C9").PasteSpecial(-4122)

# --- Column widths ---
$ws2.Columns.Item(1).ColumnWidth = 19.333333333333336
$ws2.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 26.333333333333336
$ws2.Columns.Item(4).ColumnWidth = 23.0

# --- Sheet1 selection moves to A8, and sheet2 becomes the active/selected tab ---
$ws1.Range("A8").Select()
$ws2.Range("C4").Select()
$ws2.Activate()

# --- Sheet1 gets an explicit portrait page setup ---
$ws1.PageSetup.Orientation = 1
